$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 6583.9375
$ws.Range("I28").Value = 8620.333000000001
$ws.Range("J28").Value = 474.75
$ws.Range("K28").Value = 8620.333000000001
$ws.Range("L28").Value = 474.75
$ws.Range("M28").Value = -8135.333000000001
$ws.Range("N28").Value = -1444.75

$ws.Range("H42").Value = 497
$ws.Range("I42").Value = 514.8
$ws.Range("J42").Value = 474.75
$ws.Range("K42").Value = 1544.4
$ws.Range("L42").Value = 1424.25
$ws.Range("M42").Value = -1314.4
$ws.Range("N42").Value = -1884.25

$ws.Range("H111").Value = 2252.5
$ws.Range("I111").Value = 2005.25
$ws.Range("K111").Value = 6015.75
$ws.Range("M111").Value = -2948.75

$ws.Range("H125").Value = 733.4286
$ws.Range("J125").Value = 539
$ws.Range("L125").Value = 4851
$ws.Range("N125").Value = -9771

$ws.Range("H137").Value = 2123.5833
$ws.Range("I137").Value = 1313.5
$ws.Range("J137").Value = 2933.6667
$ws.Range("K137").Value = 3940.5
$ws.Range("L137").Value = 8801.000100000001
$ws.Range("M137").Value = -1390.5
$ws.Range("N137").Value = -13901.0001

$ws.Range("H138").Value = 1955.01
$ws.Range("J138").Value = 2163.341
$ws.Range("L138").Value = 6490.022999999999
$ws.Range("N138").Value = -16770.023

$ws.Range("H141").Value = 11856.8
$ws.Range("I141").Value = 12729.777
$ws.Range("J141").Value = 4000
$ws.Range("K141").Value = 38189.331
$ws.Range("L141").Value = 12000
$ws.Range("M141").Value = -33009.331
$ws.Range("N141").Value = -22360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1430.1765
$ws.Range("I2").Value = 708.2222
$ws.Range("J2").Value = 2242.375
$ws.Range("K2").Value = 708.2222
$ws.Range("L2").Value = 2242.375
$ws.Range("M2").Value = -595.2222
$ws.Range("N2").Value = -2468.375

$ws.Range("H61").Value = 1032.6842
$ws.Range("I61").Value = 738.6
$ws.Range("K61").Value = 738.6
$ws.Range("M61").Value = -526.6

$ws.Range("H116").Value = 1430.1765
$ws.Range("I116").Value = 708.2222
$ws.Range("J116").Value = 2242.375
$ws.Range("K116").Value = 708.2222
$ws.Range("L116").Value = 2242.375
$ws.Range("M116").Value = 1585.7778
$ws.Range("N116").Value = -6830.375

$ws.Range("H136").Value = 1032.6842
$ws.Range("I136").Value = 738.6
$ws.Range("K136").Value = 2215.8
$ws.Range("M136").Value = 334.1999999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1430.1765
$ws.Range("I3").Value = 708.2222
$ws.Range("J3").Value = 2242.375
$ws.Range("K3").Value = 708.2222
$ws.Range("L3").Value = 2242.375
$ws.Range("M3").Value = -594.2222
$ws.Range("N3").Value = -2470.375

$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1617.3125
$ws.Range("I31").Value = 1274.5
$ws.Range("J31").Value = 2188.6667
$ws.Range("K31").Value = 1274.5
$ws.Range("L31").Value = 2188.6667
$ws.Range("M31").Value = -979.5
$ws.Range("N31").Value = -2778.6667

$ws.Range("H34").Value = 1617.3125
$ws.Range("I34").Value = 1274.5
$ws.Range("J34").Value = 2188.6667
$ws.Range("K34").Value = 1274.5
$ws.Range("L34").Value = 2188.6667
$ws.Range("M34").Value = -1072.5
$ws.Range("N34").Value = -2592.6667

$ws.Range("H124").Value = 10000
$ws.Range("J124").Value = 10000
$ws.Range("L124").Value = 10000
$ws.Range("N124").Value = -14910

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 143100.42
$ws.Range("I8").Value = 143100.42
$ws.Range("K8").Value = 429301.26
$ws.Range("M8").Value = -429162.26

$ws.Range("H97").Value = 512
$ws.Range("I97").Value = 512
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1536
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1040
$ws.Range("N97").ClearContents()

$ws.Range("H131").Value = 17858508
$ws.Range("I131").Value = 142857970
$ws.Range("J131").Value = 1443.9183
$ws.Range("K131").Value = 428573910
$ws.Range("L131").Value = 4331.7549
$ws.Range("M131").Value = -428568870
$ws.Range("N131").Value = -14411.7549

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3475.3572
$ws.Range("I80").Value = 2281.25
$ws.Range("J80").Value = 5067.5
$ws.Range("K80").Value = 2281.25
$ws.Range("L80").Value = 5067.5
$ws.Range("M80").Value = -1283.25
$ws.Range("N80").Value = -7063.5

$ws.Range("H83").Value = 3475.3572
$ws.Range("I83").Value = 2281.25
$ws.Range("J83").Value = 5067.5
$ws.Range("K83").Value = 11406.25
$ws.Range("L83").Value = 25337.5
$ws.Range("M83").Value = -6414.25
$ws.Range("N83").Value = -35321.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2334.1667
$ws.Range("I7").Value = 2000
$ws.Range("J7").Value = 3002.5
$ws.Range("K7").Value = 2000
$ws.Range("L7").Value = 3002.5
$ws.Range("M7").Value = -1888
$ws.Range("N7").Value = -3226.5

$ws.Range("H40").Value = 3127.7144
$ws.Range("I40").Value = 2777.8
$ws.Range("K40").Value = 2777.8
$ws.Range("M40").Value = -2641.8

$ws.Range("H68").Value = 1337.6923
$ws.Range("I68").Value = 1069.2
$ws.Range("J68").Value = 2232.6667
$ws.Range("K68").Value = 1069.2
$ws.Range("L68").Value = 2232.6667
$ws.Range("M68").Value = -320.2
$ws.Range("N68").Value = -3730.6667

$ws.Range("H71").Value = 1337.6923
$ws.Range("I71").Value = 1069.2
$ws.Range("J71").Value = 2232.6667
$ws.Range("K71").Value = 5346
$ws.Range("L71").Value = 11163.3335
$ws.Range("M71").Value = -1602
$ws.Range("N71").Value = -18651.3335

$ws.Range("H126").Value = 2334.1667
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 3002.5
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 9007.5
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -13947.5

$ws.Range("H141").Value = 50000
$ws.Range("J141").Value = 50000
$ws.Range("L141").Value = 50000
$ws.Range("N141").Value = -60360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2847.625
$ws.Range("J81").Value = 2899
$ws.Range("L81").Value = 5798
$ws.Range("N81").Value = -7920

$ws.Range("H84").Value = 2847.625
$ws.Range("J84").Value = 2899
$ws.Range("L84").Value = 28990
$ws.Range("N84").Value = -39598

$ws.Range("H136").Value = 795.75
$ws.Range("I136").Value = 573.26666
$ws.Range("K136").Value = 1719.79998
$ws.Range("M136").Value = 830.20002

Write-Output "Edit complete"
